$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "Rounded Rectangle 88" (id 89): shift right, narrower ---
$rr = $s.Shapes.Item(1)
$rr.Left = 34.18748092651367
$rr.Width = 156.0

# --- 2) "Straight Arrow Connector 114" (id 115): reroute ---
$c114 = $s.Shapes.Item(15)
$c114.Left = 88.18748474121094
$c114.Top = 72.00016021728516
$c114.Width = 0.0
$c114.Height = 37.44362258911133

# --- 3) "Straight Arrow Connector 117" (id 118): nudge right ---
$c117 = $s.Shapes.Item(18)
$c117.Left = 52.18748092651367

# --- 4) "Group 85" (id 86, client::scripts group): move up ---
$grp85 = $s.Shapes.Item(38)
$grp85.Top = 411.0625305175781

# --- 4b) Clean trailing run props on the "client::scripts" text ---
$rect86 = $grp85.GroupItems.Item(1)
$tr = $rect86.TextFrame.TextRange
$null = $tr.Delete()
$tr2 = $rect86.TextFrame.TextRange
$inserted = $tr2.InsertAfter("client::scripts")
$inserted.Font.Size = 14
$inserted.Font.Bold = -1
$inserted.Font.Color.RGB = 5287936

# --- 5) bent connector (id 162): shorten ---
$c162 = $s.Shapes.Item(39)
$c162.Height = 69.0625228881836

# --- 6) "Group 168" (id 169): move down ---
$grp168 = $s.Shapes.Item(40)
$grp168.Top = 455.0625305175781

# --- 7) "Straight Arrow Connector 161" (id 176): nudge right ---
$c161b = $s.Shapes.Item(42)
$c161b.Left = 22.187480926513672

# --- 8) New connector "Straight Arrow Connector 77", matching style of
#        the existing green dotted connectors (duplicate a non-flipped one
#        so geometry/style serialize identically, then move into place) ---
$styleSrc = $s.Shapes.Item(31)
$newConn = $styleSrc.Duplicate()
$newConn.Name = "Straight Arrow Connector 77"
$newConn.Left = 76.18748474121094
$newConn.Top = 443.23687744140625
$newConn.Width = 0.0
$newConn.Height = 17.65134048461914

# --- 9) Register the custom tag part referenced from p:custDataLst ---
$p.Tags.Add("ARTICULATE_PROJECT_OPEN", "0")
